$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(49, 8).Value = 673
$ws.Cells.Item(49, 9).Value = 509.5
$ws.Cells.Item(49, 10).Value = 1000
$ws.Cells.Item(49, 11).Value = 1528.5
$ws.Cells.Item(49, 12).Value = 3000
$ws.Cells.Item(49, 13).Value = -1392.5
$ws.Cells.Item(49, 14).Value = -3272

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(127, 8).Value = 1477.1904
$ws.Cells.Item(127, 9).Value = 499
$ws.Cells.Item(127, 10).Value = 1640.2222
$ws.Cells.Item(127, 11).Value = 1497
$ws.Cells.Item(127, 12).Value = 4920.6666
$ws.Cells.Item(127, 13).Value = 3463
$ws.Cells.Item(127, 14).Value = -14840.6666

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 2443.6128
$ws.Cells.Item(137, 9).Value = 1437.1052
$ws.Cells.Item(137, 10).Value = 4037.25
$ws.Cells.Item(137, 11).Value = 4311.3156
$ws.Cells.Item(137, 12).Value = 12111.75
$ws.Cells.Item(137, 13).Value = -1761.3156
$ws.Cells.Item(137, 14).Value = -17211.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 1315.7368
$ws.Cells.Item(74, 9).Value = 1262.8
$ws.Cells.Item(74, 10).Value = 1514.25
$ws.Cells.Item(74, 11).Value = 1262.8
$ws.Cells.Item(74, 12).Value = 1514.25
$ws.Cells.Item(74, 13).Value = -388.8
$ws.Cells.Item(74, 14).Value = -3262.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 1315.7368
$ws.Cells.Item(77, 9).Value = 1262.8
$ws.Cells.Item(77, 10).Value = 1514.25
$ws.Cells.Item(77, 11).Value = 6314
$ws.Cells.Item(77, 12).Value = 7571.25
$ws.Cells.Item(77, 13).Value = -1946
$ws.Cells.Item(77, 14).Value = -16307.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(110, 8).Value = 2869.6155
$ws.Cells.Item(110, 9).Value = 2661.6
$ws.Cells.Item(110, 10).Value = 2999.625
$ws.Cells.Item(110, 11).Value = 2661.6
$ws.Cells.Item(110, 12).Value = 2999.625
$ws.Cells.Item(110, 13).Value = -616.5999999999999
$ws.Cells.Item(110, 14).Value = -7089.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1832.375
$ws.Cells.Item(20, 9).Value = 1455
$ws.Cells.Item(20, 10).Value = 2209.75
$ws.Cells.Item(20, 11).Value = 1455
$ws.Cells.Item(20, 12).Value = 2209.75
$ws.Cells.Item(20, 13).Value = -1208
$ws.Cells.Item(20, 14).Value = -2703.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(44, 8).Value = 20000
$ws.Cells.Item(44, 9).Value = 0
$ws.Cells.Item(44, 10).Value = 20000
$ws.Cells.Item(44, 11).Value = 0
$ws.Cells.Item(44, 12).Value = 20000
$ws.Cells.Item(44, 14).Value = -20994

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(64, 8).Value = 415.83334
$ws.Cells.Item(64, 9).Value = 294
$ws.Cells.Item(64, 10).Value = 476.75
$ws.Cells.Item(64, 11).Value = 294
$ws.Cells.Item(64, 12).Value = 476.75
$ws.Cells.Item(64, 13).Value = -69
$ws.Cells.Item(64, 14).Value = -926.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(67, 8).Value = 415.83334
$ws.Cells.Item(67, 9).Value = 294
$ws.Cells.Item(67, 10).Value = 476.75
$ws.Cells.Item(67, 11).Value = 294
$ws.Cells.Item(67, 12).Value = 476.75
$ws.Cells.Item(67, 13).Value = 486
$ws.Cells.Item(67, 14).Value = -2036.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 3057.647
$ws.Cells.Item(107, 9).Value = 2671.4285
$ws.Cells.Item(107, 10).Value = 4860
$ws.Cells.Item(107, 11).Value = 2671.4285
$ws.Cells.Item(107, 12).Value = 4860
$ws.Cells.Item(107, 13).Value = -751.4285
$ws.Cells.Item(107, 14).Value = -8700

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1833.75
$ws.Cells.Item(31, 9).Value = 1414.6097
$ws.Cells.Item(31, 10).Value = 3396
$ws.Cells.Item(31, 11).Value = 1414.6097
$ws.Cells.Item(31, 12).Value = 3396
$ws.Cells.Item(31, 13).Value = -1119.6097
$ws.Cells.Item(31, 14).Value = -3986

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 1833.75
$ws.Cells.Item(34, 9).Value = 1414.6097
$ws.Cells.Item(34, 10).Value = 3396
$ws.Cells.Item(34, 11).Value = 1414.6097
$ws.Cells.Item(34, 12).Value = 3396
$ws.Cells.Item(34, 13).Value = -1212.6097
$ws.Cells.Item(34, 14).Value = -3800

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(50, 8).Value = 16283.333
$ws.Cells.Item(50, 9).Value = 8000
$ws.Cells.Item(50, 10).Value = 17940
$ws.Cells.Item(50, 11).Value = 8000
$ws.Cells.Item(50, 12).Value = 17940
$ws.Cells.Item(50, 13).Value = -7375
$ws.Cells.Item(50, 14).Value = -19190

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1058.8889
$ws.Cells.Item(58, 9).Value = 840.4815
$ws.Cells.Item(58, 10).Value = 1714.1111
$ws.Cells.Item(58, 11).Value = 840.4815
$ws.Cells.Item(58, 12).Value = 1714.1111
$ws.Cells.Item(58, 13).Value = -637.4815
$ws.Cells.Item(58, 14).Value = -2120.1111

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 2344.4285
$ws.Cells.Item(99, 9).Value = 1935
$ws.Cells.Item(99, 10).Value = 2890.3333
$ws.Cells.Item(99, 11).Value = 1935
$ws.Cells.Item(99, 12).Value = 2890.3333
$ws.Cells.Item(99, 13).Value = -437
$ws.Cells.Item(99, 14).Value = -5886.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value = 2344.4285
$ws.Cells.Item(126, 9).Value = 1935
$ws.Cells.Item(126, 10).Value = 2890.3333
$ws.Cells.Item(126, 11).Value = 5805
$ws.Cells.Item(126, 12).Value = 8670.999899999999
$ws.Cells.Item(126, 13).Value = -3335
$ws.Cells.Item(126, 14).Value = -13610.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 1058.8889
$ws.Cells.Item(136, 9).Value = 840.4815
$ws.Cells.Item(136, 10).Value = 1714.1111
$ws.Cells.Item(136, 11).Value = 2521.4445
$ws.Cells.Item(136, 12).Value = 5142.3333
$ws.Cells.Item(136, 13).Value = 28.55549999999994
$ws.Cells.Item(136, 14).Value = -10242.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(129, 8).Value = 6899.2
$ws.Cells.Item(129, 9).Value = 2086.5833
$ws.Cells.Item(129, 10).Value = 14118.125
$ws.Cells.Item(129, 11).Value = 6259.749899999999
$ws.Cells.Item(129, 12).Value = 42354.375
$ws.Cells.Item(129, 13).Value = -1259.749899999999
$ws.Cells.Item(129, 14).Value = -52354.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(137, 8).Value = 2859.077
$ws.Cells.Item(137, 9).Value = 2378.9092
$ws.Cells.Item(137, 10).Value = 5500
$ws.Cells.Item(137, 11).Value = 7136.7276
$ws.Cells.Item(137, 12).Value = 16500
$ws.Cells.Item(137, 13).Value = -2036.7276
$ws.Cells.Item(137, 14).Value = -26700

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(26, 8).Value = 13200
$ws.Cells.Item(26, 9).Value = 0
$ws.Cells.Item(26, 10).Value = 13200
$ws.Cells.Item(26, 11).Value = 0
$ws.Cells.Item(26, 12).Value = 13200
$ws.Cells.Item(26, 14).Value = -13760

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(50, 8).Value = 13200
$ws.Cells.Item(50, 9).Value = 0
$ws.Cells.Item(50, 10).Value = 13200
$ws.Cells.Item(50, 11).Value = 0
$ws.Cells.Item(50, 12).Value = 13200
$ws.Cells.Item(50, 14).Value = -14196

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2263.25
$ws.Cells.Item(102, 9).Value = 2224.5
$ws.Cells.Item(102, 10).Value = 2457
$ws.Cells.Item(102, 11).Value = 2224.5
$ws.Cells.Item(102, 12).Value = 2457
$ws.Cells.Item(102, 13).Value = -602.5
$ws.Cells.Item(102, 14).Value = -5701

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 7144280
$ws.Cells.Item(122, 9).Value = 8334493
$ws.Cells.Item(122, 10).Value = 3000
$ws.Cells.Item(122, 11).Value = 25003479
$ws.Cells.Item(122, 12).Value = 9000
$ws.Cells.Item(122, 13).Value = -25001029
$ws.Cells.Item(122, 14).Value = -13900

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3500.125
$ws.Cells.Item(7, 9).Value = 3344.7
$ws.Cells.Item(7, 10).Value = 3759.1667
$ws.Cells.Item(7, 11).Value = 3344.7
$ws.Cells.Item(7, 12).Value = 3759.1667
$ws.Cells.Item(7, 13).Value = -3232.7
$ws.Cells.Item(7, 14).Value = -3983.1667

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(32, 8).Value = 852
$ws.Cells.Item(32, 9).Value = 766.6667
$ws.Cells.Item(32, 10).Value = 980
$ws.Cells.Item(32, 11).Value = 766.6667
$ws.Cells.Item(32, 12).Value = 980
$ws.Cells.Item(32, 13).Value = -449.6667
$ws.Cells.Item(32, 14).Value = -1614

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 6690.8335
$ws.Cells.Item(40, 9).Value = 6199.2856
$ws.Cells.Item(40, 10).Value = 7379
$ws.Cells.Item(40, 11).Value = 6199.2856
$ws.Cells.Item(40, 12).Value = 7379
$ws.Cells.Item(40, 13).Value = -6063.2856
$ws.Cells.Item(40, 14).Value = -7651

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(126, 8).Value = 3500.125
$ws.Cells.Item(126, 9).Value = 3344.7
$ws.Cells.Item(126, 10).Value = 3759.1667
$ws.Cells.Item(126, 11).Value = 10034.1
$ws.Cells.Item(126, 12).Value = 11277.5001
$ws.Cells.Item(126, 13).Value = -7564.099999999999
$ws.Cells.Item(126, 14).Value = -16217.5001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(58, 8).Value = 0
$ws.Cells.Item(58, 9).Value = 0
$ws.Cells.Item(58, 10).Value = 0
$ws.Cells.Item(58, 11).Value = 0
$ws.Cells.Item(58, 12).Value = 0
$ws.Cells.Item(58, 14).ClearContents()
